# Update "想去人数" (F column) counts across the four sheets to reflect
# freshly scraped numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1397
$ws.Range("F5").Value = 5905
$ws.Range("F9").Value = 3443
$ws.Range("F10").Value = 6672
$ws.Range("F12").Value = 1319
$ws.Range("F13").Value = 765
$ws.Range("F36").Value = 17
$ws.Range("F39").Value = 1170

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 27
$ws.Range("F27").Value = 44
$ws.Range("F34").Value = 83

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 281
$ws.Range("F8").Value = 1063

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1397
$ws.Range("F10").Value = 281
$ws.Range("F11").Value = 281
$ws.Range("F14").Value = 5905
$ws.Range("F17").Value = 3443
$ws.Range("F19").Value = 6672
$ws.Range("F21").Value = 1319
$ws.Range("F24").Value = 765
$ws.Range("F26").Value = 1063
$ws.Range("F45").Value = 83
